# Elimina EC anteriores y se agregan nuevos, se modifica base de datos
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# ---------------------------------------------------------------
# 0) Preserve the two row "templates" (cell styles/borders) that
#    already live in the sheet before we start overwriting rows.
#      - row 16 : normal interior-row style
#      - row 30 : the special "closing" row style (thicker bottom
#                 border) that used to sit on the last data row
#    We copy them to two scratch rows far below the used range so
#    the originals can be freely overwritten afterwards.
# ---------------------------------------------------------------
$ws.Range("B16:J16").Copy($ws.Range("B200:J200"))
$ws.Range("B30:J30").Copy($ws.Range("B201:J201"))

# Preserve the footer rows (35 & 36) formatting/text before they are
# overwritten by new table rows.
$ws.Range("B35:J36").Copy($ws.Range("B203:J204"))

# ---------------------------------------------------------------
# 1) Header summary cells
# ---------------------------------------------------------------
$ws.Cells.Item(11,5).Value = 693034          # E11 - VALOR MORA total
$ws.Cells.Item(13,3).Value = 16              # C13 - Cant. Trabajadores
$ws.Cells.Item(13,6).Value = 6               # F13 - Cant. Periodos

# ---------------------------------------------------------------
# 2) Rebuild the worker detail table (rows 16-35)
# ---------------------------------------------------------------
# Columns: B=Tipo Doc, C=N Doc, D=Nombre, E=Periodo, F=Valor Mora, G=Salario Basico
$data = @(
    @("CC","45513862","ALBANIS ORTIZ ACEVEDO","1912",33125,828116),
    @("CC","45487604","ROSA BALVINA GONZALEZ FERIA","1912",33125,828116),
    @("CC","30878368","IRMA ESTHER POLO ARELLANO","1912",16560,392000),
    @("CC","45485491","SORELY RAMIREZ IBARRA","1912",33125,828116),
    @("CC","22802828","JACKELINE GOMEZ POLO","1912",33125,828116),
    @("CC","34970490","ENITH MARINA OSORIO GONZALEZ","1912",33125,828116),
    @("CC","51700279","NILFA DEL CARMEN RIVERA VILLERO","1912",33125,828116),
    @("CC","45468999","MARIA DE LOS REYES LUNA MUENTES","1912",15600,392000),
    @("CC","33118747","ANA ESTELA DE AVILA GAMARRA","1912",33125,828116),
    @("CC","64519634","ANA ENITH RODRIGUEZ CONTRERAS","1912",33125,828116),
    @("CC","33168342","OLGA ELENA IBAÑEZ SOLAR","1912",33125,828116),
    @("CC","22785627","BALBI ISABEL ARRIETA ANAYA","1912",15600,392000),
    @("CC","64541511","LUZ MARINA NARVAEZ MERCADO","1912",31249,392000),
    @("CC","23084112","GREGORIA DE LAS MERCEDES GUERRERO DE PUELLO","2507",56940,1423500),
    @("CC","23084112","GREGORIA DE LAS MERCEDES GUERRERO DE PUELLO","2506",56940,1423500),
    @("CC","23084112","GREGORIA DE LAS MERCEDES GUERRERO DE PUELLO","2505",56940,1423500),
    @("CC","23084112","GREGORIA DE LAS MERCEDES GUERRERO DE PUELLO","2504",56940,1423500),
    @("CC","23084112","GREGORIA DE LAS MERCEDES GUERRERO DE PUELLO","2503",56940,1423500),
    @("CC","22790768","CLARISA RAMOS PAUTT","1912",15600,392000),
    @("CC","22785865","MARGARITA ROSA JULIO LADEN","1912",15600,392000)
)

$startRow = 16
$lastRow = $startRow + $data.Length - 1     # 35

for ($i = 0; $i -lt $data.Length; $i++) {
    $r = $startRow + $i
    $row = $data[$i]

    if ($r -eq $lastRow) {
        $ws.Range("B201:J201").Copy($ws.Range("B$r`:J$r"))
    } else {
        $ws.Range("B200:J200").Copy($ws.Range("B$r`:J$r"))
    }

    $ws.Cells.Item($r,2).Value = $row[0]
    $ws.Cells.Item($r,3).Value = $row[1]
    $ws.Cells.Item($r,4).Value = $row[2]
    $ws.Cells.Item($r,5).Value = $row[3]
    $ws.Cells.Item($r,6).Value = $row[4]
    $ws.Cells.Item($r,7).Value = $row[5]
}

# ---------------------------------------------------------------
# 3) Move the signature footer from rows 35/36 down to rows 40/41
#    (rows 36-39 stay blank). The text/format was preserved in the
#    scratch rows 203/204 before the table rewrite above.
# ---------------------------------------------------------------
$ws.Range("B36:J36").ClearContents()

$ws.Range("B203:J203").Copy($ws.Range("B40:J40"))
$ws.Range("B204:J204").Copy($ws.Range("B41:J41"))

# ---------------------------------------------------------------
# 4) Column D needs to be wide enough for the new longest name
# ---------------------------------------------------------------
$ws.Columns.Item(4).ColumnWidth = 50.18

# ---------------------------------------------------------------
# 5) Clean up the scratch rows used as formatting templates
# ---------------------------------------------------------------
$ws.Range("B200:J204").Clear()

Write-Host "Edit complete"
